$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new line entries ("line7", "line8") are inserted into the name
# sequence right after "line6". Rows keep their existing row numbers
# (2-15), but since "line7"/"line8" now occupy the name slots previously
# used by "extr1"/"extr2", every row from 8 down to 15 effectively shifts
# to the next name in the sequence, and two brand new rows (16, 17) are
# appended at the bottom for "extr7"/"extr8" with their own C/D/E values.

$ws.Range("B8").Value = "line7"
$ws.Range("B9").Value = "line8"
$ws.Range("B10").Value = "extr1"
$ws.Range("B11").Value = "extr2"
$ws.Range("B12").Value = "extr3"
$ws.Range("B13").Value = "extr4"
$ws.Range("B14").Value = "extr5"
$ws.Range("B15").Value = "extr6"

# Update numeric / boolean values for rows 8-15 to match new data
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $true

$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $false

$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true

# New rows for extr7 / extr8
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false

# Copy the style used on column A (A2:A15) down to the new A16:A17 cells
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122) # xlPasteFormats
